$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update quantity for MF3RS2 Speed Cube row (B3): 3 -> 2
$ws.Range("B3").Value = 2

# Add hyperlinks to the "Link" column (E2:E5), using the URL text already
# present in each cell as the hyperlink address (the cell's existing text
# is kept as-is for display).
$ws.Hyperlinks.Add($ws.Range("E2"), $ws.Range("E2").Value2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), $ws.Range("E3").Value2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"), $ws.Range("E4").Value2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("E5"), $ws.Range("E5").Value2) | Out-Null

# Update the selected cell / active cell on the sheet
$ws.Range("E2").Select()

$wb.Save()
